# Update the dSF column (column F) with the re-pulled data values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = 5
    5  = -1
    6  = 10
    7  = 2
    8  = -3
    9  = 1
    10 = -5
    11 = 2
    12 = 4
    14 = -1
    15 = 1
    16 = 1
    17 = -3
    18 = 5
    19 = -3
    20 = 2
    21 = -1
    24 = 1
    25 = 2
    26 = 1
    27 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
